# Swap the Steps / Expected Results content between the TC3 and TC4
# test case blocks on the sheet (the TC3/TC4 id labels stay put).
#
# TC3 block: row 26 -> B26 (Steps) / D26 (Expected Results)
# TC4 block: row 33 -> B33 (Steps) / D33 (Expected Results)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tc3Steps = $ws.Range("B26").Value2
$tc3Expected = $ws.Range("D26").Value2

$tc4Steps = $ws.Range("B33").Value2
$tc4Expected = $ws.Range("D33").Value2

$ws.Range("B26").Value = $tc4Steps
$ws.Range("D26").Value = $tc4Expected

$ws.Range("B33").Value = $tc3Steps
$ws.Range("D33").Value = $tc3Expected
